$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- BOM content fix: U2 (row 39) now uses MIC94092YC6-TR instead of MIC94093YC6-TR ---
# (prevents U2 smoking per commit message)
$ws.Range("E39").Value = "MIC94092YC6-TR"
$ws.Range("J39").Value = "MIC94092_SC-70"
$ws.Range("L39").Value = "MIC94092_SC-70"

# --- BOM content fix: U3 (row 40) Value/Library Ref updated to full part number FT2232HL ---
$ws.Range("J40").Value = "FT2232HL"
$ws.Range("L40").Value = "FT2232HL"

# --- Rebuild hyperlinks: U2's digikey link now points at the new MIC94092YC6-TR part page,
#     and the old MIC94093YC6-TR hyperlink relationship is dropped ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("I24"), "http://www.digikey.com/product-detail/en/DMP3099L-7/DMP3099L-7DICT-ND/5218217")
$ws.Hyperlinks.Add($ws.Range("I25"), "http://www.digikey.com/product-detail/en/DMMT3906-7-F/DMMT3906-7-FDICT-ND/2242785")
$ws.Hyperlinks.Add($ws.Range("I38"), "http://www.digikey.com/product-detail/en/93LC56B-I%2FST/93LC56B-I%2FST-ND/319187")
$ws.Range("I39").Value = "http://www.digikey.com/product-detail/en/microchip-technology/MIC94092YC6-TR/576-3488-1-ND/2062485"
$ws.Hyperlinks.Add($ws.Range("I40"), "http://www.digikey.com/product-detail/en/FT2232HL-REEL/768-1024-1-ND/1986057")
$ws.Hyperlinks.Add($ws.Range("I41"), "http://www.digikey.com/product-detail/en/AP7365-33WG-7/AP7365-33WG-7DICT-ND/5267107")
$ws.Hyperlinks.Add($ws.Range("I42"), "http://www.digikey.com/product-detail/en/PRTR5V0U2X,215/568-4140-1-ND/1589981")
$ws.Hyperlinks.Add($ws.Range("I44"), "http://www.digikey.com/product-detail/en/ABM8-12.000MHZ-B2-T/535-9826-1-ND/2001449")
$ws.Hyperlinks.Add($ws.Range("I18"), "http://www.digikey.com/product-detail/en/CUS520,H3F/CUS520H3FCT-ND/5114381")

# --- Selection moved as part of the editing session ---
$ws.Range("E41").Select()
